# Updated symbol list on Wed Dec 28 06:16:17 UTC 2022 with GitHub Actions
# Refreshes price ("D"), and "Hora" ("G") snapshot columns for every listed coin,
# and re-syncs rows 21-27 / 41-43 where the coin ranking reshuffled (Coin/Link/Price/Volume).
#
# Values that look numeric (e.g. "243.18") are written with a leading apostrophe so
# Excel stores them as text, matching the workbook's existing text-formatted Price/Hora columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '''243.18'
$ws.Range("G2").Value = '''6'

# Row 3
$ws.Range("D3").Value = '''23.56'
$ws.Range("G3").Value = '''6'

# Row 4
$ws.Range("D4").Value = '''5.298'
$ws.Range("G4").Value = '''6'

# Row 5
$ws.Range("D5").Value = '''0.05762'
$ws.Range("G5").Value = '''6'

# Row 6
$ws.Range("D6").Value = '''6.473'
$ws.Range("G6").Value = '''6'

# Row 7
$ws.Range("D7").Value = '''3.333'
$ws.Range("G7").Value = '''6'

# Row 8
$ws.Range("D8").Value = '''0.8112'
$ws.Range("G8").Value = '''6'

# Row 9
$ws.Range("D9").Value = '''0.8749'
$ws.Range("G9").Value = '''6'

# Row 10
$ws.Range("D10").Value = '''0.1389'
$ws.Range("G10").Value = '''6'

# Row 11
$ws.Range("D11").Value = '''0.07321'
$ws.Range("G11").Value = '''6'

# Row 12
$ws.Range("D12").Value = '''0.03088'
$ws.Range("G12").Value = '''6'

# Row 13
$ws.Range("D13").Value = '''0.03057'
$ws.Range("G13").Value = '''6'

# Row 14
$ws.Range("D14").Value = '''0.09322'
$ws.Range("G14").Value = '''6'

# Row 15
$ws.Range("D15").Value = '''3.856'
$ws.Range("G15").Value = '''6'

# Row 16
$ws.Range("D16").Value = '''0.001566'
$ws.Range("G16").Value = '''6'

# Row 17
$ws.Range("D17").Value = '''0.04710'
$ws.Range("G17").Value = '''6'

# Row 18
$ws.Range("D18").Value = '''0.0006009'
$ws.Range("G18").Value = '''6'

# Row 19
$ws.Range("D19").Value = '''0.005980'
$ws.Range("G19").Value = '''6'

# Row 20
$ws.Range("D20").Value = '''0.001293'
$ws.Range("G20").Value = '''6'

# Row 21: HotbitToken
$ws.Range("B21").Value = 'HotbitToken'
$ws.Range("C21").Value = 'https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb'
$ws.Range("D21").Value = '''0.004604'
$ws.Range("E21").Value = '20HotbitTokenHTB'
$ws.Range("G21").Value = '''6'

# Row 22: NitroEx
$ws.Range("B22").Value = 'NitroEx'
$ws.Range("C22").Value = 'https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx'
$ws.Range("D22").Value = '''0.00008813'
$ws.Range("E22").Value = '21NitroExNTXBestin24h'
$ws.Range("G22").Value = '''6'

# Row 23: LEO
$ws.Range("B23").Value = 'LEO'
$ws.Range("C23").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D23").Value = '''3.578'
$ws.Range("E23").Value = '22LEOLEO'
$ws.Range("G23").Value = '''6'

# Row 24: BTSEToken
$ws.Range("B24").Value = 'BTSEToken'
$ws.Range("C24").Value = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
$ws.Range("D24").Value = '''2.141'
$ws.Range("E24").Value = '23BTSETokenBTSE'
$ws.Range("G24").Value = '''6'

# Row 25: BitpandaEcosystemToken
$ws.Range("B25").Value = 'BitpandaEcosystemToken'
$ws.Range("C25").Value = 'https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best'
$ws.Range("D25").Value = '''0.3181'
$ws.Range("E25").Value = '24BitpandaEcosystemTokenBEST'
$ws.Range("G25").Value = '''6'

# Row 26: ProBitToken
$ws.Range("B26").Value = 'ProBitToken'
$ws.Range("C26").Value = 'https://coinranking.com/coin/lQP4d6T2+probittoken-prob'
$ws.Range("D26").Value = '''0.1320'
$ws.Range("E26").Value = '25ProBitTokenPROB'
$ws.Range("G26").Value = '''6'

# Row 27: AAXToken
$ws.Range("B27").Value = 'AAXToken'
$ws.Range("C27").Value = 'https://coinranking.com/coin/LNePqkIhk+aaxtoken-aab'
$ws.Range("D27").Value = '''0.2000'
$ws.Range("E27").Value = '26AAXTokenAAB'
$ws.Range("G27").Value = '''6'

# Row 28
$ws.Range("D28").Value = '''0.0002353'
$ws.Range("G28").Value = '''6'

# Row 29
$ws.Range("G29").Value = '''6'

# Row 30
$ws.Range("G30").Value = '''6'

# Row 31
$ws.Range("G31").Value = '''6'

# Row 32
$ws.Range("G32").Value = '''6'

# Row 33
$ws.Range("G33").Value = '''6'

# Row 34
$ws.Range("G34").Value = '''6'

# Row 35
$ws.Range("G35").Value = '''6'

# Row 36
$ws.Range("G36").Value = '''6'

# Row 37
$ws.Range("G37").Value = '''6'

# Row 38
$ws.Range("G38").Value = '''6'

# Row 39
$ws.Range("G39").Value = '''6'

# Row 40
$ws.Range("D40").Value = '''0.03762'
$ws.Range("G40").Value = '''6'

# Row 41: BKEXToken
$ws.Range("B41").Value = 'BKEXToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk'
$ws.Range("D41").Value = '''0.1054'
$ws.Range("E41").Value = '40BKEXTokenBKK'
$ws.Range("G41").Value = '''6'

# Row 42: CEJI
$ws.Range("B42").Value = 'CEJI'
$ws.Range("C42").Value = 'https://coinranking.com/coin/SbKjCVJCh+ceji-ceji'
$ws.Range("D42").Value = '''0.002614'
$ws.Range("E42").Value = '41CEJICEJI'
$ws.Range("G42").Value = '''6'

# Row 43: KickToken
$ws.Range("B43").Value = 'KickToken'
$ws.Range("C43").Value = 'https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick'
$ws.Range("D43").Value = '''0.003167'
$ws.Range("E43").Value = '42KickTokenKICKWorstin24h'
$ws.Range("G43").Value = '''6'

# Row 44
$ws.Range("D44").Value = '''0.007111'
$ws.Range("G44").Value = '''6'

# Row 45
$ws.Range("D45").Value = '''0.00005487'
$ws.Range("G45").Value = '''6'

# Row 46
$ws.Range("D46").Value = '''0.00000000751'
$ws.Range("G46").Value = '''6'

# Row 47
$ws.Range("D47").Value = '''0.6009'
$ws.Range("G47").Value = '''6'

# Row 48
$ws.Range("D48").Value = '''0.001852'
$ws.Range("G48").Value = '''6'

# Row 49
$ws.Range("D49").Value = '''0.00002103'
$ws.Range("G49").Value = '''6'

# Row 50
$ws.Range("D50").Value = '''0.0002003'
$ws.Range("G50").Value = '''6'

# Row 51
$ws.Range("G51").Value = '''6'
